$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reservoirs")
$ws.Activate()

# --- U62: corrected flow magnitude ---
$ws.Range("U62").Value = 1447499

# --- Row 63: new month row (2022-10-01) ---
$ws.Range("A63").Value = 44835
$ws.Range("B63:C63").Style = "Normal"
$ws.Range("B63:T63").Value = "NaN"
$ws.Range("U63").Value = 1447499
$ws.Range("V63").Value = 554752

# --- Row 64: new month row (2022-11-01) ---
$ws.Range("A64").Value = 44866
$ws.Range("B64:C64").Style = "Normal"
$ws.Range("B64:T64").Value = "NaN"
$ws.Range("U64").Value = 1409811
$ws.Range("V64").Value = 566285

# --- Row 65: new month row (2022-12-01) ---
$ws.Range("A65").Value = 44896
$ws.Range("B65:C65").Style = "Normal"
$ws.Range("B65:T65").Value = "NaN"
$ws.Range("U65").Value = 1409811
$ws.Range("V65").Value = 586291

# --- Window view: scroll the frozen pane and update the active selection ---
$ws.Range("N47").Select()
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("W56").Select()
